$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("aiType"), shifting equip* columns right by one.
$ws.Columns.Item(3).Insert()

# The insert copies the left-neighbour's (B) formatting onto the new column;
# strip it back to the default/no style to match the data columns around it.
$ws.Range("C1:C4").ClearFormats()

# New header for column C.
$ws.Range("C1").Value = "aiType"

# Row 2: set aiType and re-affirm shifted equip values (D2:O2 already shifted by Insert()).
$ws.Range("C2").Value = 1

# Row 3: new data across C3:O3.
$row3 = @(2,1,0,0,0,2,0,1,0,2,0,2,0)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 3 + $i).Value = $row3[$i]
}

# Row 4: new data across C4:O4.
$row4 = @(0,2,0,0,0,0,0,1,0,1,0,2,0)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4[$i]
}

$ws.Range("C4").Select()
